# Auto-applied update of cryptos list (Price + Volume(1h) columns)
# Mirrors the "Updated cryptos list ... with GitHub Actions" automated refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.326.27"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.22%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.870.27"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.34%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.80%  "

$ws.Range("E6").Value = "  +0.14%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4707"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.59%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2869"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.43%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06569"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.47%  "

$ws.Range("E10").Value = "  -2.81%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08024"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.52%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "96.99"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.44%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.863.88"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.06%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.114"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.98%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6841"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.36%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "269.42"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.53%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.294.05"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.11%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.00"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.39%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007625"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.109.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.14%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.289"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.81%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.216"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.80%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.414"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.62%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.83"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.35%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.26%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.948"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.66%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.370"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.78%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09928"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.60%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.365"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.58%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.462"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.07%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.081"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.52%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04692"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.16%  "

$ws.Range("E35").Value = "  +0.27%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7005"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.28%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.704"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.04%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01878"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.62%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.635"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.75%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.307"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.14%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "71.89"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.48%  "

$ws.Range("E42").Value = "  +0.43%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8421"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.87%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4171"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.30%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.001"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.14%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.81"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.51%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.227"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.49%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.057"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.49%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "915.63"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.21%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.49"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.89%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05695"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.85%  "
